$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    # Use a fresh Range/Find each call so find-position state from a prior
    # Execute() doesn't affect the next search.
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    return $f.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Substantive correction: the practical was renumbered from 09 to 10
# ("cbnst 10 practical correction").
$r = Replace-Exact "Practical 09:" "Practical 10:"
Write-Output "Renumbered heading: $r"

# Re-assert the following paragraphs' wording unchanged; doing so via
# Find/Replace coalesces the previously fragmented runs in each paragraph
# into a single run, matching how the document re-serializes them.
Replace-Exact "Date: December 12, 2021" "Date: December 12, 2021" | Out-Null
Replace-Exact "To find root of the system equation using  Guass Seidel’s  Iteration method." "To find root of the system equation using  Guass Seidel’s  Iteration method." | Out-Null
Replace-Exact "Apply iteration method in loop with immediate updation:" "Apply iteration method in loop with immediate updation:" | Out-Null

Write-Output $d.Paragraphs(8).Range.Text
